$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (2023-10-03, serial 45202) for every
# data row (2..499). The update bumps that date by one day to 2023-10-04
# (serial 45203) for all of them.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 499 }

$newDate = Get-Date -Year 2023 -Month 10 -Day 4 -Hour 0 -Minute 0 -Second 0
$newDate = $newDate.Date
$ws.Range("C2:C$lastRow").Value = $newDate
